$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(90).Insert()

$ws.Range("A90").Value = 4
$ws.Range("B90").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C90").Value = "Los Lagos"
$ws.Range("D90").Value = 45118
$ws.Range("D90").NumberFormat = $ws.Range("D91").NumberFormat
$ws.Range("E90").Value = 10
$ws.Range("F90").Value = 100112022
$ws.Range("G90").Value = "Arveja Verde"
$ws.Range("H90").Value = "Perfection"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 40
$ws.Range("K90").Value = 42000
$ws.Range("L90").Value = 42000
$ws.Range("M90").Value = 42000
$ws.Range("N90").Value = '$/malla 25 kilos'
$ws.Range("O90").Value = 'Provincia de Limarí'
$ws.Range("P90").Value = 1680
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
